$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (sheet1) updates ---
$ws1.Range("F2").Value = 211
$ws1.Range("F3").Value = 1383
$ws1.Range("F4").Value = 19476
$ws1.Range("G4").Value = 60
$ws1.Range("F5").Value = 781
$ws1.Range("F6").Value = 298
$ws1.Range("F7").Value = 1087
$ws1.Range("F8").Value = 6
$ws1.Range("F9").Value = 7371
$ws1.Range("F11").Value = 0
$ws1.Range("F12").Value = 249
$ws1.Range("F15").Value = 0
$ws1.Range("F19").Value = 0
$ws1.Range("F21").Value = 0
$ws1.Range("F22").Value = 44
$ws1.Range("F23").Value = 48
$ws1.Range("F24").Value = 56
$ws1.Range("F25").Value = 0
$ws1.Range("F26").Value = 1062
$ws1.Range("F33").Value = 148
$ws1.Range("F35").Value = 84
$ws1.Range("F36").Value = 5
$ws1.Range("F37").Value = 12459
$ws1.Range("F38").Value = 0
$ws1.Range("I38").Value = '//i0.hdslb.com/bfs/openplatform/202407/xyYs887E1720522731548.jpeg'
$ws1.Range("F39").Value = 54
$ws1.Range("F40").Value = 0
$ws1.Range("F41").Value = 52
$ws1.Range("F42").Value = 247
$ws1.Range("F43").Value = 331
$ws1.Range("F45").Value = 316

# --- Sheet "全部类型" (sheet4) updates ---
$ws4.Range("F3").Value = 1383
$ws4.Range("F4").Value = 19476
$ws4.Range("G4").Value = 60
$ws4.Range("F7").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 479
$ws4.Range("F12").Value = 249
$ws4.Range("F13").Value = 32
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 228
$ws4.Range("F20").Value = 66
$ws4.Range("F21").Value = 0
$ws4.Range("F22").Value = 44
$ws4.Range("F23").Value = 48
$ws4.Range("F25").Value = 302
$ws4.Range("F28").Value = 6
$ws4.Range("F29").Value = 161
$ws4.Range("F30").Value = 0
$ws4.Range("F31").Value = 551
$ws4.Range("F32").Value = 2
$ws4.Range("F33").Value = 0
$ws4.Range("F35").Value = 0
$ws4.Range("F37").Value = 0
$ws4.Range("F38").Value = 0
$ws4.Range("F39").Value = 12459
$ws4.Range("F40").Value = 1313
$ws4.Range("I40").Value = '//i0.hdslb.com/bfs/openplatform/202407/xyYs887E1720522731548.jpeg'
$ws4.Range("F41").Value = 0
$ws4.Range("F42").Value = 8
$ws4.Range("F44").Value = 0
$ws4.Range("F45").Value = 331
$ws4.Range("F47").Value = 0
